$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.083.46'
$ws.Range('E2').Value = '  +4.60%  '
$ws.Range('D3').Value = '3.131.67'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'242.21"
$ws.Range('E5').Value = '  +2.55%  '
$ws.Range('D6').Value = "'612.10"
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').Value = "'1.11"
$ws.Range('E7').Value = '  +2.38%  '
$ws.Range('D8').Value = "'0.385"
$ws.Range('E8').Value = '  -1.31%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '3.129.93'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').Value = "'0.782"
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').Value = '96.734.48'
$ws.Range('E13').Value = '  +4.50%  '
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = "'34.07"
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('B16').Value = 'Toncoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D16').Value = "'5.45"
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').Value = '3.711.59'
$ws.Range('D18').Value = '3.117.63'
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('B19').Value = 'SuiNetwork'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D19').Value = "'3.52"
$ws.Range('E19').Value = '  -7.59%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = "'523.69"
$ws.Range('E20').Value = '  +19.20%  '
$ws.Range('D21').Value = "'14.60"
$ws.Range('E21').Value = '  +0.91%  '
$ws.Range('D22').Value = "'5.69"
$ws.Range('E22').Value = '  -2.17%  '
$ws.Range('D23').Value = "'0.0000193"
$ws.Range('E23').Value = '  -5.03%  '
$ws.Range('D24').Value = "'8.84"
$ws.Range('E24').Value = '  -2.98%  '
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').Value = "'5.48"
$ws.Range('E25').Value = '  -2.05%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = "'88.47"
$ws.Range('E26').Value = '  +3.47%  '
$ws.Range('D27').Value = "'11.57"
$ws.Range('E27').Value = '  -2.06%  '
$ws.Range('D28').Value = '3.297.52'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').Value = "'0.238"
$ws.Range('E30').Value = '  +2.71%  '
$ws.Range('E31').Value = '  -4.94%  '
$ws.Range('D32').Value = "'0.126"
$ws.Range('E32').Value = '  -0.81%  '
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').Value = "'9.00"
$ws.Range('E34').Value = '  -1.86%  '
$ws.Range('D35').Value = "'26.59"
$ws.Range('E35').Value = '  +3.40%  '
$ws.Range('E36').Value = '  -6.01%  '
$ws.Range('D37').Value = "'7.38"
$ws.Range('E37').Value = '  -9.03%  '
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('D39').Value = "'482.42"
$ws.Range('D40').Value = "'24.23"
$ws.Range('E40').Value = '  +1.26%  '
$ws.Range('E41').Value = '  +2.57%  '
$ws.Range('E42').Value = '  -5.01%  '
$ws.Range('E43').Value = '  -10.75%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').Value = "'3.20"
$ws.Range('E45').Value = '  -4.50%  '
$ws.Range('D46').Value = "'161.12"
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('E47').Value = '  +4.89%  '
$ws.Range('E48').Value = '  +2.45%  '
$ws.Range('D49').Value = "'4.47"
$ws.Range('D50').Value = "'44.39"
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('E51').Value = '  +0.00%  '

# Reset style on cells that were force-typed as text via the apostrophe prefix,
# so no extra quote-prefix / number-format style gets introduced.
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').Style = "Normal"
